$wb = $excel.ActiveWorkbook

# Each entry: worksheet name, cell reference, new value
# Source: "Add data for 2023-05-31" update to 2023 (column J) figures
# (and a few related prior-year corrections in columns B-I caused by reclassification)
$changes = @(
    @{ Sheet = 'Citywide Totals'; Cell = 'J2'; Value = 2931 }
    @{ Sheet = 'Citywide Totals'; Cell = 'J3'; Value = 3035 }
    @{ Sheet = 'Citywide Totals'; Cell = 'B4'; Value = 1676 }
    @{ Sheet = 'Citywide Totals'; Cell = 'C4'; Value = 1823 }
    @{ Sheet = 'Citywide Totals'; Cell = 'D4'; Value = 1951 }
    @{ Sheet = 'Citywide Totals'; Cell = 'E4'; Value = 1989 }
    @{ Sheet = 'Citywide Totals'; Cell = 'I4'; Value = 1759 }
    @{ Sheet = 'Citywide Totals'; Cell = 'J4'; Value = 683 }
    @{ Sheet = 'Citywide Totals'; Cell = 'J5'; Value = 234 }
    @{ Sheet = 'Citywide Totals'; Cell = 'J6'; Value = 3674 }
    @{ Sheet = 'Citywide Totals'; Cell = 'B7'; Value = 23308 }
    @{ Sheet = 'Citywide Totals'; Cell = 'C7'; Value = 28366 }
    @{ Sheet = 'Citywide Totals'; Cell = 'D7'; Value = 28141 }
    @{ Sheet = 'Citywide Totals'; Cell = 'E7'; Value = 25993 }
    @{ Sheet = 'Citywide Totals'; Cell = 'I7'; Value = 26206 }
    @{ Sheet = 'Citywide Totals'; Cell = 'J7'; Value = 10557 }
    @{ Sheet = 'Uptown'; Cell = 'J3'; Value = 30 }
    @{ Sheet = 'Uptown'; Cell = 'J7'; Value = 115 }
    @{ Sheet = 'Bridgeport'; Cell = 'J2'; Value = 12 }
    @{ Sheet = 'Bridgeport'; Cell = 'J7'; Value = 39 }
    @{ Sheet = 'Woodlawn'; Cell = 'J3'; Value = 50 }
    @{ Sheet = 'Woodlawn'; Cell = 'J7'; Value = 153 }
    @{ Sheet = 'North Lawndale'; Cell = 'J3'; Value = 156 }
    @{ Sheet = 'North Lawndale'; Cell = 'J7'; Value = 379 }
    @{ Sheet = 'South Deering'; Cell = 'J2'; Value = 30 }
    @{ Sheet = 'South Deering'; Cell = 'J7'; Value = 96 }
    @{ Sheet = 'New City'; Cell = 'J3'; Value = 81 }
    @{ Sheet = 'New City'; Cell = 'J6'; Value = 99 }
    @{ Sheet = 'New City'; Cell = 'J7'; Value = 280 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J2'; Value = 82 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J5'; Value = 29 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J6'; Value = 101 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J7'; Value = 320 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J8'; Value = 675 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J9'; Value = 63 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J11'; Value = 150 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J14'; Value = 39 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J16'; Value = 30 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J18'; Value = 110 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J19'; Value = 333 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J20'; Value = 215 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I29'; Value = 1555 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J29'; Value = 600 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J33'; Value = 444 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J41'; Value = 72 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J42'; Value = 423 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J43'; Value = 94 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J44'; Value = 80 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J45'; Value = 13 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J46'; Value = 37 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J48'; Value = 105 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J50'; Value = 61 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J51'; Value = 142 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J54'; Value = 208 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J57'; Value = 49 }
    @{ Sheet = 'By Neighborhood'; Cell = 'B63'; Value = 380 }
    @{ Sheet = 'By Neighborhood'; Cell = 'C63'; Value = 254 }
    @{ Sheet = 'By Neighborhood'; Cell = 'D63'; Value = 333 }
    @{ Sheet = 'By Neighborhood'; Cell = 'E63'; Value = 333 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J63'; Value = 48 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J65'; Value = 280 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J67'; Value = 379 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J72'; Value = 38 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J75'; Value = 32 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J76'; Value = 152 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J77'; Value = 94 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J79'; Value = 311 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J83'; Value = 245 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J84'; Value = 96 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J86'; Value = 62 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J88'; Value = 109 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J89'; Value = 115 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J90'; Value = 119 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J91'; Value = 123 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J99'; Value = 153 }
    @{ Sheet = 'By Neighborhood'; Cell = 'B101'; Value = 23308 }
    @{ Sheet = 'By Neighborhood'; Cell = 'C101'; Value = 28366 }
    @{ Sheet = 'By Neighborhood'; Cell = 'D101'; Value = 28141 }
    @{ Sheet = 'By Neighborhood'; Cell = 'E101'; Value = 25993 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I101'; Value = 26206 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J101'; Value = 10557 }
    @{ Sheet = 'South Chicago'; Cell = 'J3'; Value = 88 }
    @{ Sheet = 'South Chicago'; Cell = 'J7'; Value = 245 }
    @{ Sheet = 'Garfield Park'; Cell = 'J2'; Value = 117 }
    @{ Sheet = 'Garfield Park'; Cell = 'J3'; Value = 138 }
    @{ Sheet = 'Garfield Park'; Cell = 'J7'; Value = 444 }
    @{ Sheet = 'Loop'; Cell = 'J3'; Value = 41 }
    @{ Sheet = 'Loop'; Cell = 'J6'; Value = 95 }
    @{ Sheet = 'Loop'; Cell = 'J7'; Value = 208 }
    @{ Sheet = 'Englewood'; Cell = 'J2'; Value = 176 }
    @{ Sheet = 'Englewood'; Cell = 'J3'; Value = 204 }
    @{ Sheet = 'Englewood'; Cell = 'I4'; Value = 83 }
    @{ Sheet = 'Englewood'; Cell = 'J6'; Value = 158 }
    @{ Sheet = 'Englewood'; Cell = 'I7'; Value = 1555 }
    @{ Sheet = 'Englewood'; Cell = 'J7'; Value = 600 }
    @{ Sheet = 'Chatham'; Cell = 'J3'; Value = 93 }
    @{ Sheet = 'Chatham'; Cell = 'J7'; Value = 333 }
    @{ Sheet = 'Irving Park'; Cell = 'J6'; Value = 23 }
    @{ Sheet = 'Irving Park'; Cell = 'J7'; Value = 80 }
    @{ Sheet = 'Lake View'; Cell = 'J2'; Value = 20 }
    @{ Sheet = 'Lake View'; Cell = 'J7'; Value = 105 }
    @{ Sheet = 'River North'; Cell = 'J2'; Value = 23 }
    @{ Sheet = 'River North'; Cell = 'J3'; Value = 32 }
    @{ Sheet = 'River North'; Cell = 'J4'; Value = 15 }
    @{ Sheet = 'River North'; Cell = 'J7'; Value = 152 }
    @{ Sheet = 'Ashburn'; Cell = 'J3'; Value = 30 }
    @{ Sheet = 'Ashburn'; Cell = 'J4'; Value = 4 }
    @{ Sheet = 'Ashburn'; Cell = 'J7'; Value = 101 }
    @{ Sheet = 'Hermosa'; Cell = 'J3'; Value = 12 }
    @{ Sheet = 'Hermosa'; Cell = 'J7'; Value = 72 }
    @{ Sheet = 'Humboldt Park'; Cell = 'J2'; Value = 85 }
    @{ Sheet = 'Humboldt Park'; Cell = 'J5'; Value = 10 }
    @{ Sheet = 'Humboldt Park'; Cell = 'J6'; Value = 217 }
    @{ Sheet = 'Humboldt Park'; Cell = 'J7'; Value = 423 }
    @{ Sheet = 'Jefferson Park'; Cell = 'J3'; Value = 8 }
    @{ Sheet = 'Jefferson Park'; Cell = 'J7'; Value = 37 }
    @{ Sheet = 'Washington Park'; Cell = 'J3'; Value = 56 }
    @{ Sheet = 'Washington Park'; Cell = 'J7'; Value = 123 }
    @{ Sheet = 'Roseland'; Cell = 'J6'; Value = 87 }
    @{ Sheet = 'Roseland'; Cell = 'J7'; Value = 311 }
    @{ Sheet = 'Chicago Lawn'; Cell = 'J3'; Value = 64 }
    @{ Sheet = 'Chicago Lawn'; Cell = 'J7'; Value = 215 }
    @{ Sheet = 'Calumet Heights'; Cell = 'J6'; Value = 58 }
    @{ Sheet = 'Calumet Heights'; Cell = 'J7'; Value = 110 }
    @{ Sheet = 'Lincoln Square'; Cell = 'J2'; Value = 17 }
    @{ Sheet = 'Lincoln Square'; Cell = 'J7'; Value = 61 }
    @{ Sheet = 'Belmont Cragin'; Cell = 'J2'; Value = 55 }
    @{ Sheet = 'Belmont Cragin'; Cell = 'J6'; Value = 54 }
    @{ Sheet = 'Belmont Cragin'; Cell = 'J7'; Value = 150 }
    @{ Sheet = 'Avalon Park'; Cell = 'J2'; Value = 16 }
    @{ Sheet = 'Avalon Park'; Cell = 'J7'; Value = 63 }
    @{ Sheet = 'Albany Park'; Cell = 'J2'; Value = 26 }
    @{ Sheet = 'Albany Park'; Cell = 'J7'; Value = 82 }
    @{ Sheet = 'United Center'; Cell = 'J2'; Value = 25 }
    @{ Sheet = 'United Center'; Cell = 'J3'; Value = 34 }
    @{ Sheet = 'United Center'; Cell = 'J7'; Value = 109 }
    @{ Sheet = 'Austin'; Cell = 'J2'; Value = 204 }
    @{ Sheet = 'Austin'; Cell = 'J3'; Value = 215 }
    @{ Sheet = 'Austin'; Cell = 'J6'; Value = 201 }
    @{ Sheet = 'Austin'; Cell = 'J7'; Value = 675 }
    @{ Sheet = 'Armour Square'; Cell = 'J6'; Value = 19 }
    @{ Sheet = 'Armour Square'; Cell = 'J7'; Value = 29 }
    @{ Sheet = 'Streeterville'; Cell = 'J4'; Value = 29 }
    @{ Sheet = 'Streeterville'; Cell = 'J7'; Value = 62 }
    @{ Sheet = 'Pullman'; Cell = 'J2'; Value = 17 }
    @{ Sheet = 'Pullman'; Cell = 'J7'; Value = 32 }
    @{ Sheet = 'Washington Heights'; Cell = 'J2'; Value = 40 }
    @{ Sheet = 'Washington Heights'; Cell = 'J7'; Value = 119 }
    @{ Sheet = 'Little Italy, UIC'; Cell = 'J6'; Value = 44 }
    @{ Sheet = 'Little Italy, UIC'; Cell = 'J7'; Value = 142 }
    @{ Sheet = 'Mckinley Park'; Cell = 'J3'; Value = 12 }
    @{ Sheet = 'Mckinley Park'; Cell = 'J7'; Value = 49 }
    @{ Sheet = 'Hyde Park'; Cell = 'J6'; Value = 59 }
    @{ Sheet = 'Hyde Park'; Cell = 'J7'; Value = 94 }
    @{ Sheet = 'Old Town'; Cell = 'J6'; Value = 9 }
    @{ Sheet = 'Old Town'; Cell = 'J7'; Value = 38 }
    @{ Sheet = 'Riverdale'; Cell = 'J3'; Value = 32 }
    @{ Sheet = 'Riverdale'; Cell = 'J7'; Value = 94 }
    @{ Sheet = 'Jackson Park'; Cell = 'J2'; Value = 3 }
    @{ Sheet = 'Jackson Park'; Cell = 'J7'; Value = 13 }
    @{ Sheet = 'Auburn Gresham'; Cell = 'J2'; Value = 107 }
    @{ Sheet = 'Auburn Gresham'; Cell = 'J3'; Value = 91 }
    @{ Sheet = 'Auburn Gresham'; Cell = 'J6'; Value = 108 }
    @{ Sheet = 'Auburn Gresham'; Cell = 'J7'; Value = 320 }
    @{ Sheet = 'Bucktown'; Cell = 'J3'; Value = 5 }
    @{ Sheet = 'Bucktown'; Cell = 'J7'; Value = 30 }
)

foreach ($change in $changes) {
    $ws = $wb.Worksheets.Item($change.Sheet)
    $ws.Range($change.Cell).Value = $change.Value
}

Write-Output "Applied $($changes.Count) cell updates"
